$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared string text updates (volume number, date range) ---
$ws.Range("A8").Value = "Volume 31   Number  49"
$ws.Range("C9").Value = "Report Covering the Week  12/2/2024  Through  12/8/2024"

# --- Data table numeric updates ---
$ws.Range("N14").Value = -90
$ws.Range("F15").Value = 1
$ws.Range("H15").Value = -50
$ws.Range("M15").Value = 31.25
$ws.Range("C16").Value = 2
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = -33.333333333333
$ws.Range("F16").Value = 10
$ws.Range("G16").Value = 16
$ws.Range("H16").Value = -37.5
$ws.Range("I16").Value = 220
$ws.Range("J16").Value = 243
$ws.Range("K16").Value = -9.465020576131
$ws.Range("L16").Value = 20.218579234972
$ws.Range("M16").Value = 19.565217391304
$ws.Range("N16").Value = -77.059436913451
$ws.Range("C17").Value = 3
$ws.Range("E17").Value = -25
$ws.Range("F17").Value = 16
$ws.Range("G17").Value = 12
$ws.Range("H17").Value = 33.333333333333
$ws.Range("I17").Value = 245
$ws.Range("J17").Value = 225
$ws.Range("K17").Value = 8.888888888888
$ws.Range("L17").Value = 18.357487922705
$ws.Range("M17").Value = 100.819672131148
$ws.Range("N17").Value = -18.060200668896
$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 60
$ws.Range("F18").Value = 19
$ws.Range("G18").Value = 29
$ws.Range("H18").Value = -34.482758620689
$ws.Range("I18").Value = 232
$ws.Range("J18").Value = 216
$ws.Range("K18").Value = 7.407407407407
$ws.Range("L18").Value = 42.331288343558
$ws.Range("M18").Value = -3.734439834024
$ws.Range("N18").Value = -84.196185286103
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 19
$ws.Range("E19").Value = -63.157894736842
$ws.Range("F19").Value = 38
$ws.Range("G19").Value = 61
$ws.Range("H19").Value = -37.704918032786
$ws.Range("I19").Value = 671
$ws.Range("J19").Value = 704
$ws.Range("K19").Value = -4.6875
$ws.Range("L19").Value = 3.709428129829
$ws.Range("M19").Value = 53.546910755148
$ws.Range("N19").Value = -22.784810126582
$ws.Range("C20").Value = 3
$ws.Range("D20").Value = 5
$ws.Range("E20").Value = -40
$ws.Range("F20").Value = 16
$ws.Range("H20").Value = -20
$ws.Range("I20").Value = 240
$ws.Range("J20").Value = 290
$ws.Range("K20").Value = -17.241379310344
$ws.Range("L20").Value = 8.597285067873
$ws.Range("M20").Value = 25.654450261780
$ws.Range("N20").Value = -87.841945288753
$ws.Range("C21").Value = 23
$ws.Range("D21").Value = 36
$ws.Range("E21").Value = -36.111111111111
$ws.Range("F21").Value = 100
$ws.Range("G21").Value = 140
$ws.Range("H21").Value = -28.571428571428
$ws.Range("I21").Value = 1630
$ws.Range("J21").Value = 1702
$ws.Range("K21").Value = -4.230317273795
$ws.Range("L21").Value = 13.115891741845
$ws.Range("M21").Value = 36.630343671416
$ws.Range("N21").Value = -70.846002504024
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 4
$ws.Range("H22").Value = -42.857142857142
$ws.Range("I22").Value = 58
$ws.Range("J22").Value = 82
$ws.Range("K22").Value = -29.268292682926
$ws.Range("L22").Value = -23.684210526315
$ws.Range("M22").Value = 23.404255319148
$ws.Range("C24").Value = 63
$ws.Range("D24").Value = 39
$ws.Range("E24").Value = 61.538461538461
$ws.Range("F24").Value = 187
$ws.Range("G24").Value = 183
$ws.Range("H24").Value = 2.185792349726
$ws.Range("I24").Value = 1956
$ws.Range("J24").Value = 1961
$ws.Range("K24").Value = -0.254971953085
$ws.Range("L24").Value = 34.525447042641
$ws.Range("M24").Value = 120.767494356659
$ws.Range("C25").Value = 47
$ws.Range("D25").Value = 30
$ws.Range("E25").Value = 56.666666666666
$ws.Range("F25").Value = 129
$ws.Range("G25").Value = 113
$ws.Range("H25").Value = 14.159292035398
$ws.Range("I25").Value = 1327
$ws.Range("J25").Value = 1217
$ws.Range("K25").Value = 9.038619556285
$ws.Range("L25").Value = 124.915254237288
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 13
$ws.Range("E26").Value = -15.384615384615
$ws.Range("F26").Value = 45
$ws.Range("G26").Value = 44
$ws.Range("H26").Value = 2.272727272727
$ws.Range("I26").Value = 530
$ws.Range("J26").Value = 492
$ws.Range("K26").Value = 7.723577235772
$ws.Range("L26").Value = 0.760456273764
$ws.Range("M26").Value = 11.814345991561
$ws.Range("F27").Value = 1
$ws.Range("H27").Value = -50
$ws.Range("C28").Value = 1
$ws.Range("C28").NumberFormat = "#,##0"
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = -50
$ws.Range("I28").Value = 61
$ws.Range("J28").Value = 79
$ws.Range("K28").Value = -22.784810126582
$ws.Range("L28").Value = -32.967032967033
